$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The existing row 21 (previous weekly entry) is pushed down to row 22
# unchanged, and row 21 is updated with a new weekly entry's values
# (this is a weekly data refresh, per the commit message).

# Copy current row 21 values down to new row 22 first (before overwriting row 21).
$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R")
foreach ($col in $cols) {
    $ws.Range("$col" + "22").Value2 = $ws.Range("$col" + "21").Value2
}
$ws.Range("D22").NumberFormat = $ws.Range("D21").NumberFormat

# Now update row 21 with the new weekly entry's values.
$ws.Range("D21").Value2 = 44504
$ws.Range("J21").Value2 = 100
$ws.Range("K21").Value2 = 17000
$ws.Range("L21").Value2 = 18000
$ws.Range("M21").Value2 = 17500
$ws.Range("O21").Value2 = "Región del Maule"
$ws.Range("P21").Value2 = 700
